$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.029551440430632
$ws.Range("D2").Value = 1.032645670186132
$ws.Range("E2").Value = 1.039604559829757
$ws.Range("F2").Value = 1.051184082651919
$ws.Range("I2").Value = 1.033831422049696
$ws.Range("J2").Value = 1.034697588555688
$ws.Range("K2").Value = 1.035450302334172
$ws.Range("L2").Value = 1.042389274677821
$ws.Range("M2").Value = 1.053936289174165
$ws.Range("N2").Value = 1.03616697781291

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.030510146849641
$ws.Range("D3").Value = 1.033335809851813
$ws.Range("E3").Value = 1.040483903499396
$ws.Range("F3").Value = 1.052214318521658
$ws.Range("I3").Value = 1.034015289360558
$ws.Range("J3").Value = 1.035297327688827
$ws.Range("K3").Value = 1.035949574867607
$ws.Range("L3").Value = 1.043078688937177
$ws.Range("M3").Value = 1.054778560566304
$ws.Range("N3").Value = 1.036767568644409

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.031130712055624
$ws.Range("D4").Value = 1.033782176940289
$ws.Range("E4").Value = 1.041053466987981
$ws.Range("F4").Value = 1.052881588695312
$ws.Range("I4").Value = 1.034132518772976
$ws.Range("J4").Value = 1.035685019294703
$ws.Range("K4").Value = 1.036271763254257
$ws.Range("L4").Value = 1.043524718803655
$ws.Range("M4").Value = 1.055323598082971
$ws.Range("N4").Value = 1.037155810816834

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031391649130359
$ws.Range("D5").Value = 1.033969780623271
$ws.Range("E5").Value = 1.041293046670854
$ws.Range("F5").Value = 1.053162260936616
$ws.Range("I5").Value = 1.034181383615067
$ws.Range("J5").Value = 1.035847913108804
$ws.Range("K5").Value = 1.036407000992311
$ws.Range("L5").Value = 1.043712212763193
$ws.Range("M5").Value = 1.055552738439669
$ws.Range("N5").Value = 1.037318935958835

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031435464609703
$ws.Range("D6").Value = 1.034001277193481
$ws.Range("E6").Value = 1.041333281049314
$ws.Range("F6").Value = 1.053209395915862
$ws.Range("I6").Value = 1.034189563683093
$ws.Range("J6").Value = 1.035875258293221
$ws.Range("K6").Value = 1.036429695646351
$ws.Range("L6").Value = 1.043743692801303
$ws.Range("M6").Value = 1.055591212471389
$ws.Range("N6").Value = 1.037346319976551

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031134198508301
$ws.Range("D7").Value = 1.033784683904225
$ws.Range("E7").Value = 1.041056667733611
$ws.Range("F7").Value = 1.052885338456139
$ws.Range("I7").Value = 1.034133173352232
$ws.Range("J7").Value = 1.035687196251334
$ws.Range("K7").Value = 1.036273571134904
$ws.Range("L7").Value = 1.043527224175031
$ws.Range("M7").Value = 1.055326659842996
$ws.Range("N7").Value = 1.037157990864993

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029875394161851
$ws.Range("D8").Value = 1.032878946630497
$ws.Range("E8").Value = 1.039901619590535
$ws.Range("F8").Value = 1.051532122820806
$ws.Range("I8").Value = 1.033893921805998
$ws.Range("J8").Value = 1.034900351327207
$ws.Range("K8").Value = 1.035619214293819
$ws.Range("L8").Value = 1.042622278844136
$ws.Range("M8").Value = 1.054220931601241
$ws.Range("N8").Value = 1.036370028530821

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027658922334064
$ws.Range("D9").Value = 1.031281451451806
$ws.Range("E9").Value = 1.037870686233235
$ws.Range("F9").Value = 1.049152521079001
$ws.Range("I9").Value = 1.033458993086386
$ws.Range("J9").Value = 1.033510958317053
$ws.Range("K9").Value = 1.034459498964473
$ws.Range("L9").Value = 1.041027173643551
$ws.Range("M9").Value = 1.052272782248971
$ws.Range("N9").Value = 1.034978662423227

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026182460433297
$ws.Range("D10").Value = 1.030215541666793
$ws.Range("E10").Value = 1.036519759191306
$ws.Range("F10").Value = 1.047569498971449
$ws.Range("I10").Value = 1.033160110464511
$ws.Range("J10").Value = 1.032582813400236
$ws.Range("K10").Value = 1.033681935540267
$ws.Range("L10").Value = 1.039963501754656
$ws.Range("M10").Value = 1.05097426405748
$ws.Range("N10").Value = 1.034049199434169

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025543425510594
$ws.Range("D11").Value = 1.029753789977661
$ws.Range("E11").Value = 1.035935524269749
$ws.Range("F11").Value = 1.046884846367009
$ws.Range("I11").Value = 1.033028580377324
$ws.Range("J11").Value = 1.032180479487848
$ws.Range("K11").Value = 1.033344204784484
$ws.Range("L11").Value = 1.039502866863812
$ws.Range("M11").Value = 1.05041206153644
$ws.Range("N11").Value = 1.033646294161477

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025306102549525
$ws.Range("D12").Value = 1.029582245144944
$ws.Range("E12").Value = 1.035718623395931
$ws.Range("F12").Value = 1.046630657959477
$ws.Range("I12").Value = 1.032979407506803
$ws.Range("J12").Value = 1.032030969225751
$ws.Range("K12").Value = 1.033218601041284
$ws.Range("L12").Value = 1.039331758640525
$ws.Range("M12").Value = 1.05020324501915
$ws.Range("N12").Value = 1.033496571577658

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025357007172451
$ws.Range("D13").Value = 1.029619043418305
$ws.Range("E13").Value = 1.035765144376293
$ws.Range("F13").Value = 1.046685176698132
$ws.Range("I13").Value = 1.032989969578708
$ws.Range("J13").Value = 1.032063042638344
$ws.Range("K13").Value = 1.033245550493895
$ws.Range("L13").Value = 1.039368462280765
$ws.Range("M13").Value = 1.050248036380736
$ws.Range("N13").Value = 1.033528690538176

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.025523807438835
$ws.Range("D14").Value = 1.02973961062708
$ws.Range("E14").Value = 1.035917592926304
$ws.Range("F14").Value = 1.046863832578888
$ws.Range("I14").Value = 1.033024522190701
$ws.Range("J14").Value = 1.032168122255597
$ws.Range("K14").Value = 1.033333825507234
$ws.Range("L14").Value = 1.039488723161598
$ws.Range("M14").Value = 1.050394800478227
$ws.Range("N14").Value = 1.033633919380539

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.025626584311328
$ws.Range("D15").Value = 1.029813892164753
$ws.Range("E15").Value = 1.03601153611407
$ws.Range("F15").Value = 1.046973924563623
$ws.Range("I15").Value = 1.033045769240799
$ws.Range("J15").Value = 1.032232856613814
$ws.Range("K15").Value = 1.033388194075928
$ws.Range("L15").Value = 1.039562818834177
$ws.Range("M15").Value = 1.05048522810345
$ws.Range("N15").Value = 1.033698745668969

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.026224877084013
$ws.Range("D16").Value = 1.030246182374915
$ws.Range("E16").Value = 1.036558548308203
$ws.Range("F16").Value = 1.047614954176964
$ws.Range("I16").Value = 1.033168795234538
$ws.Range("J16").Value = 1.03260950574696
$ws.Range("K16").Value = 1.033704327764881
$ws.Range("L16").Value = 1.039994071421994
$ws.Range("M16").Value = 1.051011577017349
$ws.Range("N16").Value = 1.034075929687086

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026600246456671
$ws.Range("D17").Value = 1.030517292383015
$ws.Range("E17").Value = 1.036901869527246
$ws.Range("F17").Value = 1.048017271652955
$ws.Range("I17").Value = 1.033245401156879
$ws.Range("J17").Value = 1.032845650233389
$ws.Range("K17").Value = 1.033902352058733
$ws.Range("L17").Value = 1.040264569822235
$ws.Range("M17").Value = 1.051341759635508
$ws.Range("N17").Value = 1.034312409525772

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.02681922031104
$ws.Range("D18").Value = 1.030675406334225
$ws.Range("E18").Value = 1.037102192909293
$ws.Range("F18").Value = 1.048252014377097
$ws.Range("I18").Value = 1.033289880286693
$ws.Range("J18").Value = 1.032983346714568
$ws.Range("K18").Value = 1.034017755813837
$ws.Range("L18").Value = 1.040422341319007
$ws.Range("M18").Value = 1.051534355729453
$ws.Range("N18").Value = 1.034450301551749

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026893889341296
$ws.Range("D19").Value = 1.030729315688202
$ws.Range("E19").Value = 1.037170509835362
$ws.Range("F19").Value = 1.048332068739611
$ws.Range("I19").Value = 1.033305011934102
$ws.Range("J19").Value = 1.033030290407948
$ws.Range("K19").Value = 1.034057088445628
$ws.Range("L19").Value = 1.040476136340894
$ws.Range("M19").Value = 1.051600027039091
$ws.Range("N19").Value = 1.034497311910558

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02655997005264
$ws.Range("D20").Value = 1.030488206897529
$ws.Range("E20").Value = 1.036865027169765
$ws.Range("F20").Value = 1.047974098746961
$ws.Range("I20").Value = 1.033237203143238
$ws.Range("J20").Value = 1.03282031856611
$ws.Range("K20").Value = 1.033881116311521
$ws.Range("L20").Value = 1.040235548481843
$ws.Range("M20").Value = 1.051306333522119
$ws.Range("N20").Value = 1.034287041884621

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025474687724452
$ws.Range("D21").Value = 1.029704107399211
$ws.Range("E21").Value = 1.035872697583704
$ws.Range("F21").Value = 1.046811219498383
$ws.Range("I21").Value = 1.033014356046793
$ws.Range("J21").Value = 1.032137180730557
$ws.Range("K21").Value = 1.033307835001884
$ws.Range("L21").Value = 1.039453309541774
$ws.Range("M21").Value = 1.050351581811971
$ws.Range("N21").Value = 1.033602933914983

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02479257707838
$ws.Range("D22").Value = 1.029210941707778
$ws.Range("E22").Value = 1.035249417446781
$ws.Range("F22").Value = 1.046080777771126
$ws.Range("I22").Value = 1.032872410946751
$ws.Range("J22").Value = 1.031707286064051
$ws.Range("K22").Value = 1.032946491011613
$ws.Range("L22").Value = 1.038961439469003
$ws.Range("M22").Value = 1.049751352868188
$ws.Range("N22").Value = 1.033172428748742

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025154153003958
$ws.Range("D23").Value = 1.029472393998332
$ws.Range("E23").Value = 1.035579769338137
$ws.Range("F23").Value = 1.046467931472605
$ws.Range("I23").Value = 1.032947832224832
$ws.Range("J23").Value = 1.031935217023624
$ws.Range("K23").Value = 1.033138131257872
$ws.Range("L23").Value = 1.039222193141469
$ws.Range("M23").Value = 1.050069539488618
$ws.Range("N23").Value = 1.033400683396422

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.026578169141718
$ws.Range("D24").Value = 1.030501349438126
$ws.Range("E24").Value = 1.036881674429257
$ws.Range("F24").Value = 1.047993606484619
$ws.Range("I24").Value = 1.033240908102484
$ws.Range("J24").Value = 1.032831764987827
$ws.Range("K24").Value = 1.033890712142072
$ws.Range("L24").Value = 1.040248661993828
$ws.Range("M24").Value = 1.051322341038853
$ws.Range("N24").Value = 1.034298504561569

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028231727428399
$ws.Range("D25").Value = 1.031694608680545
$ws.Range("E25").Value = 1.038395202073393
$ws.Range("F25").Value = 1.049767113945675
$ws.Range("I25").Value = 1.033573009269704
$ws.Range("J25").Value = 1.033870485083777
$ws.Range("K25").Value = 1.034760096184394
$ws.Range("L25").Value = 1.041439597648385
$ws.Range("M25").Value = 1.035338699759199
